$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 header: "Cost" -> "Cost / ticket" (E5). B5/D5 stay "Schedule"/"Duration".
$ws.Range("E5").Value = "Cost / ticket"

# Row 6: first flight option becomes the Alaska/Horizon flight (previously on row 7)
$ws.Range("A6").Value = "AlaskaOperated by Horizon Air as Alaska Horizon"
$ws.Range("B6").Value = "7:00 PM – 7:26 AM+1"
$ws.Range("D6").Value = "14 hr 26 min"
$ws.Range("E6").Value = 439

# Row 7: second flight option becomes American
$ws.Range("A7").Value = "American"
$ws.Range("B7").Value = "7:43 AM – 12:45 PM"
$ws.Range("D7").Value = "7 hr 2 min"
$ws.Range("E7").Value = 512

# Active cell selection moves from J9 to H9
$ws.Range("H9").Select()
